$wb = $excel.ActiveWorkbook
